$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A29").Value = "2024-11-04 11:01:41"
$ws.Range("B29").Value = "Success"
$ws.Range("C29").Value = @'

    <!DOCTYPE html>
    <html lang="en">
    <head>
        <meta charset="UTF-8">
         <style>
            body { font-family: Arial, sans-serif; margin: 20px; }
            h1 { color: #333; }
            h2 { color: #555; margin-top: 20px; }
            table { width: 80%; border-collapse: collapse; margin-bottom: 20px; }
            th, td { border: 1px solid #ddd; padding: 10px; text-align: left; }
            th { background-color: #f2f2f2; color: #333; }
            tr:nth-child(even) { background-color: #f9f9f9; }
            tr:hover { background-color: #ddd; }
            /* Specific styles for averages */
            .averages-table {
                width: 80%;
                margin-top: 20px;
                border: 2px solid #007BFF; /* blue border for emphasis */
            }
            .averages-table th {
                background-color: #007BFF; /* blue background for headers */
                color: white; /* white text for headers */
                font-weight: bold;
            }
            .averages-table td {
                background-color: #E7F1FF; /* light blue background for data */
            }
        </style>
    </head>
    <body>
    <h1>Weather Report</h1>
    <h2>Averages:</h2>
    <table class="averages-table">
    <tr>
    <th>Kaupunki</th><th>Keskimääräinen lämpötila</th><th>Keskimääräinen tuulen nopeus</th><th>Säätila</th></tr><tr><td>Helsinki</td><td>1.7 Celsius</td><td>14.0 m/s</td><td>Ei merkittäviä sääilmiöitä</td></tr><tr><td>Espoo</td><td>1.2 Celsius</td><td>9.3 m/s</td><td>Ei merkittäviä sääilmiöitä</td></tr><tr><td>Vantaa</td><td>0.3 Celsius</td><td>9.9 m/s</td><td>Ei merkittäviä sääilmiöitä</td></tr></table><h2>Source: OpenWeatherMap</h2><table border='1'><tr><th>Kaupunki</th><th>Lämpötila</th><th>Säätila</th><th>Tuulen_nopeus</th></tr><tr><td>Helsinki</td><td>1.12</td><td>moderate rain</td><td>8.05</td></tr><tr><td>Espoo</td><td>0.8</td><td>moderate rain</td><td>5.14</td></tr><tr><td>Vantaa</td><td>0.67</td><td>moderate rain</td><td>5.14</td></tr></table><h2>Source: Finnish Meteorological Institute</h2><table border='1'><tr><th>Kaupunki</th><th>Lämpötila</th><th>Säätila</th><th>Tuulen_nopeus</th></tr><tr><td>Helsinki</td><td>2.2 Celsius</td><td>Ei merkittäviä sääilmiöitä</td><td>5.9 m/s avarage speed measured in the last 10 minutes</td></tr><tr><td>Espoo</td><td>1.7 Celsius</td><td>Ei merkittäviä sääilmiöitä</td><td>4.2 m/s avarage speed measured in the last 10 minutes</td></tr><tr><td>Vantaa</td><td>-0.1 Celsius</td><td>Ei merkittäviä sääilmiöitä</td><td>4.8 m/s avarage speed measured in the last 10 minutes</td></tr></table></body>
'@
